{"js": "// Change 1: FirstParagraph (paragraph 3) \u2014 tighten recruitment/exclusion description.\n{\n  const oldText = \"Like the pilot study, participants were recruited on Amazon Mechanical Turk for a study on decision-making and performance, with an initial sample of 1296 before excluding participants who did not meet inclusion criteria. The inclusion criteria were nearly identical to those in the pilot study, with the exception that participants were not excluded if they failed the comprehension check questions. Thus, a total of 284 participants were excluded before analyses: 25 were excluded because they did not indicate they were American or lived in the United States, 3 were excluded for indicating \u201cOther\u201d for their gender, 192 were excluded for using a phone or tablet to complete the survey, and 64 were excluded for an incomplete survey. The final sample consisted of 1012 participants (53.66% women), with an average age of 37.66 (\";\n  const newText = \"Like the pilot study, we recruited workers on Amazon Mechanical Turk for a study on decision-making and performance. The pre-screening criteria were nearly identical to those in the pilot study, with the exception that workers were not excluded if they failed the comprehension questions to increase power. The final sample consisted of 1012 participants (53.66% women), with an average age of 37.66 (\";\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Change 1: expected 1 match, found \" + results.items.length);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Change 2: BodyText paragraph describing the multiplication task / condition assignment.\n{\n  const oldText = \"Participants who met the inclusion criteria were told they would be completing a two-minute multiplication task where they would be able to choose how they would be paid for their performance. The multiplication task consisted of multiplying two numbers with digits ranging from 1-12 (e.g., 1 X 5, 12 X 11) as quickly as possible. Then, they were provided examples and had to complete three comprehension check questions, which they had to pass to proceed. After completing the comprehension check questions, participants were assigned to either a knowledge of preparation condition or a control condition based on their gender. Participants in the knowledge of preparation condition were presented the following text:\";\n  const newText = \"Participants were told they would be completing a two-minute multiplication task where they would be able to choose how they would be paid for their performance. For the task, participants answered questions from the multiplication tables with numbers ranging from 1-12 (e.g., 1 X 5, 12 X 11) as quickly as possible. Then, they were provided examples and had to complete three comprehension questions, which they had to pass to proceed. After completing the comprehension questions, participants were randomly assigned to either a \u201cknowledge of preparation\u201d condition or a control condition based on their gender. Participants in the knowledge of preparation condition were presented the following text:\";\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Change 2: expected 1 match, found \" + results.items.length);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Change 3: BodyText paragraph describing payment schemes / practice procedure.\n{\n  const oldText = \"Participants assigned to the control condition simply proceeded without seeing this text. Then, all participants learned about the payment schemes (either piece-rate or tournament) for the multiplication task and had to pass several comprehension check questions about the payment schemes before being given the opportunity to choose a payment scheme. The payment schemes followed the same structure as the payment schemes in the pilot study, with the exception that the payment was doubled for each scheme (i.e., $.20 per problem in the tournament scheme, $.10 per problem in the piece-rate scheme). In the preparation condition, participants were reminded that they had the option to prepare before completing the task, while participants in the control condition did not have this reminder. Then, participants made a payment scheme choice, where the order of the presentation of the tournament and piece-rate payment options were randomized for each condition, so the tournament payment scheme was listed first for some participants, while the piece-rate payment scheme was listed first for others. After choosing a payment scheme, participants in both conditions were given the chance to prepare before the multiplication task. If they agreed to practice (described as the choice/decision to practice in subsequent analyses), participants were asked, for each times table, if they wanted to practice problems from that specific times table. If they chose to practice a specific times table, they had the chance to practice all twelve combinations of numbers for that times table. They could only proceed if they answered all practice questions correctly. Then, they were asked if they would like to continue practicing or move onto the next times table, while a review table was displayed. This process was repeated for each times table. The practice and review table for each times table was presented in sequential order (i.e., starting at the one times table up to the twelve times table). We measured the number of rounds of preparation each participant completed for analyses (i.e., total practice count), which was calculated as the total number of times a participant agreed to complete a round of preparation (including the choice to repeat a table and the choice to prepare in the first place). Once finished practicing, participants completed as many problems as possible from the paid multiplication task for two minutes and received feedback about their absolute (but not relative) performance.\";\n  const newText = \"Participants assigned to the control condition simply proceeded without seeing this text. Then, all participants learned about the possible payment schemes that they could choose (either piece-rate or tournament) and had to pass several comprehension questions about the payment schemes before being choosing a payment scheme. For the tournament scheme, participants were paid $.20 per problem they answered correctly only if they beat a randomly assigned partner, while the piece-rate scheme paid participants $.10 per problem, regardless of other participants\u2019 performance. In the preparation condition, participants were reminded that they had the option to prepare before completing the task, while participants in the control condition did not have this reminder. Then, participants made a payment scheme choice, where the order of presentation of the tournament and piece-rate payment options was randomized and counter-balanced for each condition. After choosing a payment scheme, participants in both conditions were given the chance to prepare before the multiplication task. If they chose to practice (described as the choice/decision to practice in subsequent analyses), participants were asked, for each multiplication table, if they wanted to practice problems from that specific multiplication table. If they chose to practice a specific multiplication table, they had the chance to practice all twelve combinations of numbers for that multiplication table. They could only proceed if they answered all practice questions correctly. Then, they were asked if they would like to continue practicing or move onto the next multiplication table, while a review table was displayed. This process was repeated for each multiplication table. The practice and review table for each multiplication table was presented in sequential order (i.e., starting at the 1 multiplication table up to the 12 multiplication table). We measured the number of rounds of preparation each participant completed for analyses (i.e., total practice count), which was calculated as the total number of times a participant agreed to complete a round of preparation (including the choice to repeat a table and the choice to prepare in the first place). Once finished practicing, participants completed as many problems as possible from the paid multiplication task for two minutes and received feedback about their absolute (but not relative) performance.\";\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Change 3: expected 1 match, found \" + results.items.length);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Change 1: tighten recruitment/exclusion description (FirstParagraph)\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = 'Like the pilot study, participants were recruited on Amazon Mechanical Turk for a study on decision-making and performance, with an initial sample of 1296 before excluding participants who did not meet inclusion criteria. The inclusion criteria were nearly identical to those in the pilot study, with the exception that participants were not excluded if they failed the comprehension check questions. Thus, a total of 284 participants were excluded before analyses: 25 were excluded because they did not indicate they were American or lived in the United States, 3 were excluded for indicating \u201cOther\u201d for their gender, 192 were excluded for using a phone or tablet to complete the survey, and 64 were excluded for an incomplete survey. The final sample consisted of 1012 participants (53.66% women), with an average age of 37.66 ('\n$find.Replacement.Text = 'Like the pilot study, we recruited workers on Amazon Mechanical Turk for a study on decision-making and performance. The pre-screening criteria were nearly identical to those in the pilot study, with the exception that workers were not excluded if they failed the comprehension questions to increase power. The final sample consisted of 1012 participants (53.66% women), with an average age of 37.66 ('\n$result1 = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\nif (-not $result1) {\n    throw \"Change 1: Find.Execute did not find/replace the target text\"\n}\n\n# Change 2: multiplication task / condition-assignment paragraph\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = 'Participants who met the inclusion criteria were told they would be completing a two-minute multiplication task where they would be able to choose how they would be paid for their performance. The multiplication task consisted of multiplying two numbers with digits ranging from 1-12 (e.g., 1 X 5, 12 X 11) as quickly as possible. Then, they were provided examples and had to complete three comprehension check questions, which they had to pass to proceed. After completing the comprehension check questions, participants were assigned to either a knowledge of preparation condition or a control condition based on their gender. Participants in the knowledge of preparation condition were presented the following text:'\n$find.Replacement.Text = 'Participants were told they would be completing a two-minute multiplication task where they would be able to choose how they would be paid for their performance. For the task, participants answered questions from the multiplication tables with numbers ranging from 1-12 (e.g., 1 X 5, 12 X 11) as quickly as possible. Then, they were provided examples and had to complete three comprehension questions, which they had to pass to proceed. After completing the comprehension questions, participants were randomly assigned to either a \u201cknowledge of preparation\u201d condition or a control condition based on their gender. Participants in the knowledge of preparation condition were presented the following text:'\n$result2 = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\nif (-not $result2) {\n    throw \"Change 2: Find.Execute did not find/replace the target text\"\n}\n\n# Change 3: payment schemes / practice procedure paragraph\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = 'Participants assigned to the control condition simply proceeded without seeing this text. Then, all participants learned about the payment schemes (either piece-rate or tournament) for the multiplication task and had to pass several comprehension check questions about the payment schemes before being given the opportunity to choose a payment scheme. The payment schemes followed the same structure as the payment schemes in the pilot study, with the exception that the payment was doubled for each scheme (i.e., $.20 per problem in the tournament scheme, $.10 per problem in the piece-rate scheme). In the preparation condition, participants were reminded that they had the option to prepare before completing the task, while participants in the control condition did not have this reminder. Then, participants made a payment scheme choice, where the order of the presentation of the tournament and piece-rate payment options were randomized for each condition, so the tournament payment scheme was listed first for some participants, while the piece-rate payment scheme was listed first for others. After choosing a payment scheme, participants in both conditions were given the chance to prepare before the multiplication task. If they agreed to practice (described as the choice/decision to practice in subsequent analyses), participants were asked, for each times table, if they wanted to practice problems from that specific times table. If they chose to practice a specific times table, they had the chance to practice all twelve combinations of numbers for that times table. They could only proceed if they answered all practice questions correctly. Then, they were asked if they would like to continue practicing or move onto the next times table, while a review table was displayed. This process was repeated for each times table. The practice and review table for each times table was presented in sequential order (i.e., starting at the one times table up to the twelve times table). We measured the number of rounds of preparation each participant completed for analyses (i.e., total practice count), which was calculated as the total number of times a participant agreed to complete a round of preparation (including the choice to repeat a table and the choice to prepare in the first place). Once finished practicing, participants completed as many problems as possible from the paid multiplication task for two minutes and received feedback about their absolute (but not relative) performance.'\n$find.Replacement.Text = 'Participants assigned to the control condition simply proceeded without seeing this text. Then, all participants learned about the possible payment schemes that they could choose (either piece-rate or tournament) and had to pass several comprehension questions about the payment schemes before being choosing a payment scheme. For the tournament scheme, participants were paid $.20 per problem they answered correctly only if they beat a randomly assigned partner, while the piece-rate scheme paid participants $.10 per problem, regardless of other participants\u2019 performance. In the preparation condition, participants were reminded that they had the option to prepare before completing the task, while participants in the control condition did not have this reminder. Then, participants made a payment scheme choice, where the order of presentation of the tournament and piece-rate payment options was randomized and counter-balanced for each condition. After choosing a payment scheme, participants in both conditions were given the chance to prepare before the multiplication task. If they chose to practice (described as the choice/decision to practice in subsequent analyses), participants were asked, for each multiplication table, if they wanted to practice problems from that specific multiplication table. If they chose to practice a specific multiplication table, they had the chance to practice all twelve combinations of numbers for that multiplication table. They could only proceed if they answered all practice questions correctly. Then, they were asked if they would like to continue practicing or move onto the next multiplication table, while a review table was displayed. This process was repeated for each multiplication table. The practice and review table for each multiplication table was presented in sequential order (i.e., starting at the 1 multiplication table up to the 12 multiplication table). We measured the number of rounds of preparation each participant completed for analyses (i.e., total practice count), which was calculated as the total number of times a participant agreed to complete a round of preparation (including the choice to repeat a table and the choice to prepare in the first place). Once finished practicing, participants completed as many problems as possible from the paid multiplication task for two minutes and received feedback about their absolute (but not relative) performance.'\n$result3 = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\nif (-not $result3) {\n    throw \"Change 3: Find.Execute did not find/replace the target text\"\n}\n"}
